# "+ now the software work + no touch"
#
# Typography sheet: the "Large" typography row no longer overrides the
# sheet-wide wildcard characters, so H5's custom value is cleared.
#
# Translation sheet: the single remaining on-screen text (row 4,
# SingleUseId1) now uses the Default typography, is left-aligned, and its
# translated text changes from "DC_Load" to "DC Load Control Panel". The
# second text (row 5, SingleUseId2 / "<> %" - the touch-related glyph) is
# removed entirely now that the software runs without a touch screen.

$wb = $excel.ActiveWorkbook
$wsTypography  = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet -------------------------------------------------
# Clear H5 (Wildcard Ranges override for the "Large" typography) while
# keeping the cell present (empty) rather than removing it outright.
$wsTypography.Range("H5").ClearContents()
$wsTypography.Range("H5").Style = "Normal"

# --- Translation sheet -------------------------------------------------
# Row 4 (TEXT ID = SingleUseId1): switch typography to Default, alignment
# to Left, and update the translated text.
$wsTranslation.Range("C4").Value = "Default"
$wsTranslation.Range("D4").Value = "Left"
$wsTranslation.Range("F4").Value = "DC Load Control Panel"

# Row 5 (TEXT ID = SingleUseId2) is dropped entirely.
$wsTranslation.Range("B5:F5").ClearContents()

Write-Host "texts.xlsx updated: cleared Typography!H5, Translation row 4 now Default/Left/'DC Load Control Panel', Translation row 5 removed"
